$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 9.5
$ws.Range("H3").Value = 5.25
$ws.Range("J3").Value = 9.5
$ws.Range("L3").Value = 1.8
